$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 801-802; existing rows 801-817 shift down to 803-819.
$ws.Range("A801:R802").EntireRow.Insert()

# Row 801 - new record (Primera)
$ws.Cells.Item(801, 1).Value = 6
$ws.Cells.Item(801, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(801, 3).Value = "Metropolitana"
$ws.Cells.Item(801, 4).Value = 45239
$ws.Cells.Item(801, 5).Value = 13
$ws.Cells.Item(801, 6).Value = 100112052
$ws.Cells.Item(801, 7).Value = "Albahaca"
$ws.Cells.Item(801, 8).Value = "Sin especificar"
$ws.Cells.Item(801, 9).Value = "Primera"
$ws.Cells.Item(801, 10).Value = 250
$ws.Cells.Item(801, 11).Value = 7000
$ws.Cells.Item(801, 12).Value = 7000
$ws.Cells.Item(801, 13).Value = 7000
$ws.Cells.Item(801, 14).Value = "$/docena de matas"
$ws.Cells.Item(801, 15).Value = "Región Metropolitana"
$ws.Cells.Item(801, 16).Value = 1167
$ws.Cells.Item(801, 17).Value = 6
$ws.Cells.Item(801, 18).Value = "Hortaliza"

# Row 802 - new record (Segunda)
$ws.Cells.Item(802, 1).Value = 6
$ws.Cells.Item(802, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(802, 3).Value = "Metropolitana"
$ws.Cells.Item(802, 4).Value = 45239
$ws.Cells.Item(802, 5).Value = 13
$ws.Cells.Item(802, 6).Value = 100112052
$ws.Cells.Item(802, 7).Value = "Albahaca"
$ws.Cells.Item(802, 8).Value = "Sin especificar"
$ws.Cells.Item(802, 9).Value = "Segunda"
$ws.Cells.Item(802, 10).Value = 150
$ws.Cells.Item(802, 11).Value = 6000
$ws.Cells.Item(802, 12).Value = 6000
$ws.Cells.Item(802, 13).Value = 6000
$ws.Cells.Item(802, 14).Value = "$/docena de matas"
$ws.Cells.Item(802, 15).Value = "Región Metropolitana"
$ws.Cells.Item(802, 16).Value = 1000
$ws.Cells.Item(802, 17).Value = 6
$ws.Cells.Item(802, 18).Value = "Hortaliza"
